# Add Rubin and Shalin ballots to the "ballots" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# ---- Row 44: Roger Rubin ----
$ws.Range("A44").Value = "Roger Rubin"
$ws.Range("C44").Value = "x"
$ws.Range("D44").Value = "x"
$ws.Range("E44").Value = "x"
$ws.Range("I44").Value = "x"
$ws.Range("O44").Value = "x"
$ws.Range("Q44").Value = "x"
$ws.Range("AK44").Value = 6
$ws.Range("AL44").Value = "Newsday+Twitter"

# Copy the date style from the row above so the new date cell keeps the
# existing date number format (instead of Excel inventing a new one).
$ws.Range("AM43").Copy()
$ws.Range("AM44").PasteSpecial(-4122)
$ws.Range("AM44").Value = 43448

# ---- Row 45: Mike Shalin ----
$ws.Range("A45").Value = "Mike Shalin"
$ws.Range("C45").Value = "x"
$ws.Range("D45").Value = "x"
$ws.Range("E45").Value = "x"
$ws.Range("H45").Value = "x"
$ws.Range("I45").Value = "x"
$ws.Range("K45").Value = "x"
$ws.Range("O45").Value = "x"
$ws.Range("Q45").Value = "x"
$ws.Range("T45").Value = "x"
$ws.Range("V45").Value = "x"
$ws.Range("AK45").Value = 10
$ws.Range("AL45").Value = "Twitter"

$ws.Range("AM43").Copy()
$ws.Range("AM45").PasteSpecial(-4122)
$ws.Range("AM45").Value = 43448

# Reflect the new last row as the active selection, matching where the
# author ended up after entering the new ballots.
$ws.Range("B45").Select() | Out-Null
